$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K5").Value = "NAT"
$ws.Range("K6").Value = "NAT"
$ws.Range("K7").Value = "NAT"

$ws.Range("K5").Select()
